$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(1.165638061119211,1.12590394751075,1.102162474457373,1.092652554004417,1.091083403256789,1.102033552432886,1.151801751470458,1.254599708302663,1.33331006234647,1.369812482673467,1.383735265837515,1.380732294115205,1.370955912495106,1.364980631499861,1.33093855598446,1.310233228774734,1.298389651558864,1.294390885480652,1.312430559537518,1.37382475639339,1.414532893899036,1.392752829034578,1.311436958801778,1.226231833635154)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colB[$i]
}

$colC = @(0.04494834593693042,0.04002713358481458,0.03702309864640085,0.03580326045799609,0.03560096596370954,0.03700663007612093,0.04324779679421908,0.05563144357428484,0.06482577282906732,0.06903111565981135,0.07062697020795383,0.07028312226924527,0.06916233918300918,0.06847627078911955,0.0645514114135608,0.06214955005272316,0.06077020233929886,0.06030354318097864,0.06240501012023003,0.06949144770723592,0.07414260902318404,0.07165835353360706,0.06228951191567944,0.05226503726890996)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $colC[$i]
}

$colE = @(0.05070288413446633,0.05018988964175719,0.04989810068810385,0.0497850366576067,0.04976661563880391,0.0498965521973247,0.05052119622875928,0.05192979600561998,0.05307639177306811,0.05362221722349503,0.05383238455114281,0.05378696685423634,0.05363943822503003,0.05354952495321541,0.05304120760520803,0.05273557120636596,0.05256205875512521,0.05250370256947079,0.05276787066430089,0.05368267672334426,0.05430080637237822,0.0539690491328777,0.05275326122012558,0.05152909770908209)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $colE[$i]
}

$colF = @(0.4443680307746263,0.387822817061874,0.3531389305168915,0.3390132514313251,0.336668177824194,0.3529483938344953,0.4248636149813336,0.5661985755041457,0.6702781546542269,0.7176906081379002,0.7356546913071611,0.7317853510981394,0.7191683204515869,0.7114413442032514,0.6671810134426437,0.6400460337125793,0.6244449056556647,0.619163680173358,0.642933953830422,0.7228739723491628,0.7751780083420101,0.7472568307830727,0.6416283278902171,0.5279251897347166)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $colF[$i]
}

$colG = @(0.002565582288970012,0.002569468254315551,0.002571980270614005,0.002573035726711025,0.002573212907294287,0.00257199437597144,0.00256689607561944,0.002557893611445314,0.002551879790931882,0.002549272918055439,0.002548304186005467,0.002548512001502958,0.002549192850786604,0.002549612288642018,0.002552052740904948,0.002553582812692701,0.002554475002516366,0.002554779169266162,0.002553418678900651,0.002548992369294885,0.002546206927272121,0.002547683773009156,0.002553492844691418,0.002560223134609817)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $colG[$i]
}

$colI = @(1.743079149791846,1.745451868291909,1.747624983234104,1.748690524126239,1.748878322945778,1.747638624941679,1.74374852240777,1.741810608963597,1.743868429022648,1.745563482994228,1.746314697063468,1.746148043805263,1.745623093460139,1.745315790690569,1.743772942031157,1.743020977914099,1.742659891027564,1.742549894183682,1.743093632328581,1.74577431520477,1.748163673237684,1.746830041127239,1.743060563399986,1.741724375130893)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item(2 + $i, 9).Value = $colI[$i]
}

$colK = @(0.9663782970679051,0.9251258962553948,0.9003606299513081,0.8904102881414246,0.8887665978009807,0.90022586251294,0.9520373686662822,1.058124315259988,1.138824952206505,1.176142722718879,1.190361511656306,1.187295353688626,1.177310759501864,1.171206287790824,1.136398365890244,1.115200296324105,1.103064819976737,1.098965758224608,1.117450957266755,1.180241104923397,1.221787354038383,1.1995666992743,1.116433272809644,1.028942348323596)
for ($i = 0; $i -lt $colK.Length; $i++) {
    $ws.Cells.Item(2 + $i, 11).Value = $colK[$i]
}

$colL = @(0.2567060454407226,0.2496110860395646,0.245401772636427,0.2437233960554721,0.2434469343091621,0.2453789878622388,0.2542291761806581,0.2727529487186899,0.2870794454794918,0.2937537843692581,0.2963038347211437,0.2957536291602878,0.2939631249942494,0.29286933647235,0.2866464260690265,0.2828691385453652,0.2807113321613031,0.2799832753283056,0.2832697064091576,0.2944884252749631,0.3019523839018348,0.2979566530769659,0.2830885666765539,0.2676161924353266)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $ws.Cells.Item(2 + $i, 12).Value = $colL[$i]
}

$colN = @(2.900895565848131,2.920152422369171,2.932765770687595,2.93810418827735,2.939002605072119,2.932836963471047,2.907371357603125,2.8637053167445,2.835456790334874,2.82344018034766,2.819009909522578,2.81995869931184,2.823073289630742,2.824996722160492,2.836258904177015,2.843381624981433,2.847556889539376,2.848984030778681,2.842615276577973,2.822155195779288,2.809483851147789,2.816182602815203,2.842961492499072,2.874845733845333)
for ($i = 0; $i -lt $colN.Length; $i++) {
    $ws.Cells.Item(2 + $i, 14).Value = $colN[$i]
}

